$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from the "Database - ydb463" paragraph;
#    the same bookmark id is re-created below, wrapping the new password
#    value instead.
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$codeRPr = '<w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Times New Roman"/><w:color w:val="D4D4D4"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr>'

# ------------------------------------------------------------------
# 2. Create one placeholder empty paragraph right after "Database - ydb463",
#    then flood it with all of the new paragraphs' OOXML in a single
#    InsertXML call (the hyperlink itself is added as a placeholder marker
#    "ZZLINKZZ" that gets turned into a real Hyperlinks.Add call below --
#    Hyperlinks.Add on a *collapsed* range snaps to the start of its
#    paragraph in this host, so it must target a real, non-collapsed
#    Find match instead).
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()

$body  = "<w:p $wNs/>"

$body += "<w:p $wNs>"
$body +=   "<w:pPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"1E1E1E`"/><w:spacing w:line=`"285`" w:lineRule=`"atLeast`"/></w:pPr>"
$body +=   "<w:r><w:t xml:space=`"preserve`">My </w:t></w:r>"
$body +=   "<w:proofErr w:type=`"spellStart`"/>"
$body +=   "<w:r><w:t>gmail</w:t></w:r>"
$body +=   "<w:proofErr w:type=`"spellEnd`"/>"
$body +=   "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
$body +=   "<w:proofErr w:type=`"spellStart`"/>"
$body +=   "<w:r><w:t>accout</w:t></w:r>"
$body +=   "<w:proofErr w:type=`"spellEnd`"/>"
$body +=   "<w:r><w:t xml:space=`"preserve`"> for sending letters &#8211;</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p $wNs>"
$body +=   "<w:pPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"1E1E1E`"/><w:spacing w:line=`"285`" w:lineRule=`"atLeast`"/>$codeRPr</w:pPr>"
$body +=   "<w:r><w:t xml:space=`"preserve`">Mail -  </w:t></w:r>"
$body +=   "<w:r><w:rPr><w:rFonts w:ascii=`"Consolas`" w:eastAsia=`"Times New Roman`" w:hAnsi=`"Consolas`" w:cs=`"Times New Roman`"/><w:color w:val=`"CE9178`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr><w:t>godofcode991@gmail.com</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p $wNs>"
$body +=   "<w:pPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"1E1E1E`"/><w:spacing w:line=`"285`" w:lineRule=`"atLeast`"/>$codeRPr</w:pPr>"
$body +=   "<w:r><w:t xml:space=`"preserve`">Password - </w:t></w:r>"
$body +=   "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>"
$body +=   "<w:r><w:rPr><w:rFonts w:ascii=`"Consolas`" w:eastAsia=`"Times New Roman`" w:hAnsi=`"Consolas`" w:cs=`"Times New Roman`"/><w:color w:val=`"6A9955`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr><w:t>Z3RynAgdSzC8iFj</w:t></w:r>"
$body +=   "<w:bookmarkEnd w:id=`"0`"/>"
$body += "</w:p>"

$body += "<w:p $wNs>"
$body +=   "<w:r><w:t xml:space=`"preserve`">Link for on\off secure app access to google account(it option must be in &#8220;on&#8221; state for giving access to email-sending for different apps - </w:t></w:r>"
$body +=   "<w:r><w:t>ZZLINKZZ</w:t></w:r>"
$body +=   "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
$body += "</w:p>"

$body += "<w:p $wNs>"
$body +=   "<w:pPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"1E1E1E`"/><w:spacing w:after=`"0`" w:line=`"285`" w:lineRule=`"atLeast`"/>$codeRPr</w:pPr>"
$body +=   "<w:r>$codeRPr<w:t xml:space=`"preserve`"> </w:t></w:r>"
$body +=   "<w:r><w:rPr><w:rFonts w:ascii=`"Consolas`" w:eastAsia=`"Times New Roman`" w:hAnsi=`"Consolas`" w:cs=`"Times New Roman`"/><w:color w:val=`"9CDCFE`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr><w:t>username:</w:t></w:r>"
$body +=   "<w:r>$codeRPr<w:t xml:space=`"preserve`"> </w:t></w:r>"
$body +=   "<w:r><w:rPr><w:rFonts w:ascii=`"Consolas`" w:eastAsia=`"Times New Roman`" w:hAnsi=`"Consolas`" w:cs=`"Times New Roman`"/><w:color w:val=`"CE9178`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr><w:t>'godofcode991@gmail.com'</w:t></w:r>"
$body +=   "<w:r>$codeRPr<w:t>,</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p $wNs>"
$body +=   "<w:pPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"1E1E1E`"/><w:spacing w:after=`"0`" w:line=`"285`" w:lineRule=`"atLeast`"/>$codeRPr</w:pPr>"
$body +=   "<w:r>$codeRPr<w:t xml:space=`"preserve`">        </w:t></w:r>"
$body += "</w:p>"

$body += "<w:p $wNs>"
$body +=   "<w:pPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"1E1E1E`"/><w:spacing w:after=`"0`" w:line=`"285`" w:lineRule=`"atLeast`"/>$codeRPr</w:pPr>"
$body +=   "<w:r>$codeRPr<w:t xml:space=`"preserve`">            </w:t></w:r>"
$body +=   "<w:r><w:rPr><w:rFonts w:ascii=`"Consolas`" w:eastAsia=`"Times New Roman`" w:hAnsi=`"Consolas`" w:cs=`"Times New Roman`"/><w:color w:val=`"9CDCFE`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr><w:t>password:</w:t></w:r>"
$body +=   "<w:r>$codeRPr<w:t xml:space=`"preserve`"> </w:t></w:r>"
$body +=   "<w:r><w:rPr><w:rFonts w:ascii=`"Consolas`" w:eastAsia=`"Times New Roman`" w:hAnsi=`"Consolas`" w:cs=`"Times New Roman`"/><w:color w:val=`"CE9178`"/><w:sz w:val=`"21`"/><w:szCs w:val=`"21`"/></w:rPr><w:t>'Z3RynAgdSzC8iFj'</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p $wNs/>"

$placeholder = $d.Paragraphs.Item(3).Range
$placeholder.InsertXML($body)

Write-Output "Inserted body, paragraph count now: $($d.Paragraphs.Count)"

# ------------------------------------------------------------------
# 3. Turn the "ZZLINKZZ" placeholder into a real hyperlink (creates the
#    external relationship too), then fix its character style to match
#    the document's existing hyperlink style ("a3" == built-in Hyperlink).
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("ZZLINKZZ")
if (-not $found) {
    Write-Output "ERROR: placeholder not found"
}
$url = "https://myaccount.google.com/lesssecureapps?pli=1&rapt=AEjHL4Nx9F51r-bJJXmc1Uq2pLAcCGMI-UNufXV1LFsI0dDyKugALUfJcIs7kaNuLpRzMZk2MU3KtK0jkVXJ7ApyXbK9ghrZ2w"
$h = $d.Hyperlinks.Add($findRange, $url, "", "", $url)
$h.Range.Style = "a3"

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
